$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# register_data: add a second registered user (santa+1@gmail.com), mirroring
# the existing row 2 sample data / formatting.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("register_data")

$ws2.Range("A3").Value = "male"
$ws2.Range("B3").Value = "santa"
$ws2.Range("C3").Value = "santhosh"
$ws2.Range("D3").Value = 22
$ws2.Range("E3").Value = "March "
$ws2.Range("F3").Value = 2003
$ws2.Range("G3").Value = "santa+1@gmail.com"
$ws2.Range("H3").Value = 123456
$ws2.Range("I3").Value = 123456

$ws2.Range("A3:I3").Font.Name = "Arial"
$ws2.Range("A3:I3").Font.Size = 10
$ws2.Range("A3:I3").Borders.LineStyle = 1

$ws2.Hyperlinks.Add($ws2.Range("G3"), "mailto:santa+1@gmail.com", "", "", "santa+1@gmail.com")
$ws2.Range("G3").Font.Color = $ws2.Range("G2").Font.Color
$ws2.Range("G3").Font.Underline = $false

# ---------------------------------------------------------------------------
# Billing_Address: fix up the "country" header label and append a blank row
# so the sheet has a spot ready for the next registrant's address.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Billing_Address")
$ws3.Range("D1").Value = "country"
$ws3.Range("A1:H1").RowHeight = 12.8

$ws3.Range("A3:H3").Value = ""

# ---------------------------------------------------------------------------
# Active sheet: register_data is now the sheet of interest.
# ---------------------------------------------------------------------------
$ws2.Activate()
